$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to be stored as Text before assigning a numeric-looking
# string value, so Excel does not silently convert it into a Number (which would
# drop formatting such as trailing zeros, e.g. "5.090" -> 5.09).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

Set-TextValue $ws.Range('D2') '27.405.72'
$ws.Range('E2').Value = '  +1.49%  '
Set-TextValue $ws.Range('D3') '1.859.56'
$ws.Range('E3').Value = '  +0.51%  '
$ws.Range('E4').Value = '  -0.03%  '
Set-TextValue $ws.Range('D5') '311.20'
$ws.Range('E5').Value = '  +0.55%  '
$ws.Range('E6').Value = '  -0.12%  '
Set-TextValue $ws.Range('D7') '0.4765'
$ws.Range('E7').Value = '  -0.30%  '
Set-TextValue $ws.Range('D8') '0.3806'
$ws.Range('E8').Value = '  +3.50%  '
Set-TextValue $ws.Range('D9') '0.07311'
$ws.Range('E9').Value = '  +1.07%  '
Set-TextValue $ws.Range('D10') '0.9301'
$ws.Range('E10').Value = '  -0.14%  '
Set-TextValue $ws.Range('D11') '20.79'
$ws.Range('E11').Value = '  +5.01%  '
Set-TextValue $ws.Range('D12') '0.07784'
$ws.Range('E12').Value = '  +0.59%  '
Set-TextValue $ws.Range('D13') '1.843.18'
$ws.Range('E13').Value = '  -0.91%  '
Set-TextValue $ws.Range('D14') '5.437'
$ws.Range('E14').Value = '  +1.57%  '
Set-TextValue $ws.Range('D15') '6.538'
$ws.Range('E15').Value = '  +1.34%  '
Set-TextValue $ws.Range('D16') '90.03'
Set-TextValue $ws.Range('D17') '1.012'
$ws.Range('E17').Value = '  -0.32%  '
Set-TextValue $ws.Range('D18') '0.000008810'
$ws.Range('E18').Value = '  +1.97%  '
$ws.Range('E19').Value = '  -0.21%  '
Set-TextValue $ws.Range('D20') '27.386.13'
$ws.Range('E20').Value = '  +1.38%  '
Set-TextValue $ws.Range('D21') '14.62'
$ws.Range('E21').Value = '  +0.53%  '
Set-TextValue $ws.Range('D22') '5.090'
$ws.Range('E22').Value = '  +0.27%  '
$ws.Range('E23').Value = '  +0.37%  '
Set-TextValue $ws.Range('D24') '1.942'
$ws.Range('E24').Value = '  +0.44%  '
Set-TextValue $ws.Range('D25') '155.59'
$ws.Range('E25').Value = '  +1.71%  '
$ws.Range('E26').Value = '  +1.27%  '
Set-TextValue $ws.Range('D27') '2.008'
$ws.Range('E27').Value = '  -0.39%  '
Set-TextValue $ws.Range('D28') '115.26'
$ws.Range('E28').Value = '  +0.80%  '
Set-TextValue $ws.Range('D29') '4.943'
$ws.Range('E29').Value = '  -0.45%  '
Set-TextValue $ws.Range('D30') '0.08895'
$ws.Range('E30').Value = '  +0.37%  '
Set-TextValue $ws.Range('D31') '3.322'
$ws.Range('E31').Value = '  +0.33%  '
$ws.Range('E32').Value = '  +1.90%  '
Set-TextValue $ws.Range('D33') '4.588'
$ws.Range('E33').Value = '  +1.75%  '
Set-TextValue $ws.Range('D34') '0.7496'
$ws.Range('E34').Value = '  +1.17%  '
Set-TextValue $ws.Range('D35') '2.716'
$ws.Range('E35').Value = '  -0.93%  '
$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D36') '0.02041'
$ws.Range('E36').Value = '  +4.12%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range('D37') '1.121'
$ws.Range('E37').Value = '  +0.96%  '
Set-TextValue $ws.Range('D38') '0.5532'
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D39') '0.05256'
$ws.Range('E39').Value = '  -0.20%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range('D40') '2.988'
$ws.Range('E40').Value = '  +0.43%  '
Set-TextValue $ws.Range('D41') '7.022'
$ws.Range('E41').Value = '  -0.08%  '
Set-TextValue $ws.Range('D42') '8.594'
$ws.Range('E43').Value = '  +0.24%  '
Set-TextValue $ws.Range('D44') '0.4866'
$ws.Range('E44').Value = '  +2.31%  '
Set-TextValue $ws.Range('D45') '10.59'
$ws.Range('E45').Value = '  -0.12%  '
$ws.Range('E46').Value = '  -0.23%  '
Set-TextValue $ws.Range('D47') '1.663'
$ws.Range('E47').Value = '  +3.41%  '
Set-TextValue $ws.Range('D48') '102.84'
$ws.Range('E48').Value = '  +1.01%  '
Set-TextValue $ws.Range('D49') '67.24'
$ws.Range('E49').Value = '  +2.10%  '
Set-TextValue $ws.Range('D50') '0.06096'
Set-TextValue $ws.Range('D51') '0.9115'
$ws.Range('E51').Value = '  +2.44%  '
